# "Refined metadata to be additional tab"
#
# 1. Refresh the "time_taken" timestamps on the existing "data" sheet
#    (cells F2:F12) with the new panel_query_time.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    panel data source that was queried.

$wb   = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. refresh F2:F12 timestamps on the "data" sheet -------------------
$timestamps = @(
    "2021-10-05 14:34:01.485445",
    "2021-10-05 14:34:01.485453",
    "2021-10-05 14:34:01.485456",
    "2021-10-05 14:34:01.485459",
    "2021-10-05 14:34:01.485461",
    "2021-10-05 14:34:01.485464",
    "2021-10-05 14:34:01.485467",
    "2021-10-05 14:34:01.485469",
    "2021-10-05 14:34:01.485472",
    "2021-10-05 14:34:01.485475",
    "2021-10-05 14:34:01.485477"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2. add the "metadata" worksheet -------------------------------------
# Duplicate "data" (right after itself) so the new sheet inherits the same
# sheet/page setup (outline + page-setup properties, margins, header
# style, …), then strip it down and overwrite with the metadata content.
$data.Copy($null, $data)
$meta = $wb.Worksheets.Item(2)
$meta.Name = "metadata"

# Drop the gene rows inherited from "data" - only two rows survive.
$meta.Range("A3:F12").Clear()

# Extend the styled header row out to column G (reuse the same header
# formatting as the rest of row 1).
$data.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Cells.Item(2, 1).Value = 0
$meta.Range("B2").Value = "Hirschsprung disease"
$meta.Range("C2").Value = 110
# keep "0.18" as text (not a number) like the source data export
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.18"
$meta.Range("E2").Value = "2021-08-03T22:26:06.454717Z"
$meta.Range("F2").Value = "2021-10-05 14:34:01.481610"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/110/?format=json"

[void]$meta.Range("A1").Select()

# restore "data" as the active sheet/tab
[void]$data.Activate()
[void]$data.Range("A1").Select()
